$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data range as Text so numeric-looking strings
# (e.g. "1.005", "307.26") are preserved verbatim instead of
# being coerced into floating point numbers on assignment.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.472.46"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.656.97"
$ws.Range("E3").Value = "  -2.87%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "307.26"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.3611"
$ws.Range("E7").Value = "  -3.61%  "
$ws.Range("D8").Value = "47.41"
$ws.Range("E8").Value = "  -3.89%  "
$ws.Range("D9").Value = "0.3247"
$ws.Range("E9").Value = "  -5.77%  "
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("D11").Value = "0.06985"
$ws.Range("E11").Value = "  -6.60%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "5.881"
$ws.Range("E13").Value = "  -5.90%  "
$ws.Range("D14").Value = "19.37"
$ws.Range("E14").Value = "  -7.56%  "
$ws.Range("D15").Value = "1.661.59"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "6.560"
$ws.Range("E16").Value = "  -5.86%  "
$ws.Range("D17").Value = "0.00001041"
$ws.Range("E17").Value = "  -7.73%  "
$ws.Range("D18").Value = "0.06535"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "76.31"
$ws.Range("E20").Value = "  -9.56%  "
$ws.Range("D21").Value = "5.913"
$ws.Range("E21").Value = "  -6.70%  "
$ws.Range("D22").Value = "15.62"
$ws.Range("E22").Value = "  -8.98%  "
$ws.Range("D23").Value = "12.53"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").Value = "24.473.72"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "2.470"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "2.315"
$ws.Range("E26").Value = "  -16.50%  "
$ws.Range("D27").Value = "146.81"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").Value = "18.43"
$ws.Range("E28").Value = "  -8.94%  "
$ws.Range("D29").Value = "1.844.94"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").Value = "123.81"
$ws.Range("E30").Value = "  -5.74%  "
$ws.Range("D31").Value = "1.170"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "3.971"
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("D33").Value = "5.613"
$ws.Range("E33").Value = "  -17.41%  "
$ws.Range("D34").Value = "1.700"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("D35").Value = "0.08410"
$ws.Range("E35").Value = "  -5.28%  "
$ws.Range("D36").Value = "12.29"
$ws.Range("E36").Value = "  -10.20%  "
$ws.Range("D37").Value = "5.179"
$ws.Range("E37").Value = "  -6.46%  "
$ws.Range("D38").Value = "0.06039"
$ws.Range("E38").Value = "  -8.38%  "
$ws.Range("D39").Value = "0.02201"
$ws.Range("E39").Value = "  -7.86%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.202"
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("D41").Value = "8.206"
$ws.Range("E41").Value = "  -8.49%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.2049"
$ws.Range("E42").Value = "  -8.24%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "0.5887"
$ws.Range("E44").Value = "  -8.72%  "
$ws.Range("D45").Value = "3.736"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("E46").Value = "  -9.36%  "
$ws.Range("D47").Value = "0.5578"
$ws.Range("E47").Value = "  -8.98%  "
$ws.Range("D48").Value = "121.95"
$ws.Range("E48").Value = "  -6.17%  "
$ws.Range("D49").Value = "1.930"
$ws.Range("E49").Value = "  -8.98%  "
$ws.Range("D50").Value = "0.06896"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").Value = "74.13"
$ws.Range("E51").Value = "  -6.57%  "

# Restore the default (unstyled) cell style now that the values
# are committed as text, so no stray style indices are introduced.
$ws.Range("B2:E51").Style = "Normal"

